# Consolidated error messages V2.3.xlsx - "Updated error message for IDA"
#
# Semantic changes (reconstructed from the OOXML diff, resolving shared-string
# renumbering noise):
#   1. IDA sheet, row 14 (Sl No 13) column B ("Scenario"): the trailing
#      "; Failure in Decryption" clause is dropped from the scenario text.
#   2. IDA sheet, row 34 (Sl No 40) column C ("Message"): the message text is
#      shortened from "Unable to decrypt Authentication Request." to
#      "Unable to decrypt Request.", and the whole row (B:F) is highlighted
#      yellow to flag the update (matching the existing highlighted row 67).
#   3. Selection/view state on the IDA sheet moves to J7 with no scrolled
#      top-left cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("IDA")

# 1) Row 14 - Scenario text shortened (drop "; Failure in Decryption")
$ws.Range("B14").Value = "Could not process request/Unknown error; Invalid Auth Request"

# 2) Row 34 - Message text shortened
$ws.Range("C34").Value = "“Unable to decrypt Request.”"

# Highlight the updated row (B34:F34) yellow, same as the other flagged row (67)
$ws.Range("B34:F34").Interior.Color = 65535

# 3) Update the active selection / scroll position on the IDA sheet
$ws.Activate()
$ws.Range("J7").Select()

$wb.Save()
